# carbon_guessr workbook update
# - Fix mislabeled source #1 ("Mazac (2022)" -> "Poore & Nemecek (2018)" style title)
# - Add a new source (Sanchez 2020, Fairphone 5 LCA report) with hyperlink
# - Add a new "Manufacturing" category data row (Fairphone 5 smartphone) to the
#   data table, expanding Table1 to match
# - Restore cursor/selection state on the touched sheets

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "sources" sheet: correct source #1's title and append source #3
# ---------------------------------------------------------------------------
$sources = $wb.Worksheets.Item("sources")

$sources.Range("B2").Value = "Poor & Nemecek (2018)"

$sources.Range("A4").Value = 3
$sources.Range("B4").Value = "Sánchez (2020)"
$sources.Range("C4").Value = 'Sánchez, D., S. J. Baur, and L. Eguren. "Life Cycle Assessment of the Fairphone 5. Berlin: Fraunhofer IZM." 19.06. 2020–Version 2 David Sánchez Sarah-Jane Baur Lara Eguren 5.3 (2020): 83.'
$sources.Range("D4").Value = "Sánchez (2020)"
$sources.Hyperlinks.Add($sources.Range("D4"), "https://www.fairphone.com/wp-content/uploads/2024/09/Fairphone5_LCA_Report_2024.pdf")
$sources.Range("D4").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# 2. "data" sheet: append a new "Manufacturing" category row for the
#    Fairphone 5 smartphone, growing the Table1 ListObject along the way
# ---------------------------------------------------------------------------
$data = $wb.Worksheets.Item("data")
$table = $data.ListObjects.Item(1)
$table.ListRows.Add() | Out-Null

$data.Range("B56").Value = "Fairphone 5 smartphone"
$data.Range("C56").Value = "Unit"
$data.Range("D56").Value = "Average carbon footprint for the production and transportation of the Fairphone 5 smartphone. This does not include the use phase and the 'end of life' phase (such as recycling and waste management). Note that Fairphones are some of the most sustainably manufactured smartphones, and their carbon footprint is significantly lower than those of other typical smartphones as a result."
$data.Range("E56").Value = 30.41
$data.Range("F56").Value = "Manufacturing"
$data.Range("G56").Value = 3

# ---------------------------------------------------------------------------
# 3. Restore the on-screen selections: "sources" was left on D5, "data" is
#    the active tab selected at B57
# ---------------------------------------------------------------------------
$sources.Activate()
$sources.Range("D5").Select()

$data.Activate()
$data.Range("B57").Select()
